# Generate Report for Handoff
# Refresh the localization-status report: the handoff package UUID/hash
# changed from 4c3ad086-a9de-4630-8da0-a2acb5a7658a (hash
# 328c4a9fa66d971c6e8210ba0abc9d999c7c06f7) to
# bf2ba36c-0e8e-4eb8-8e64-e8dc319da171 (hash
# f7f65d20bcc6f852bb855cade74d66448e1a21ee), and the handoff timestamps
# advanced. Update cell text + the matching hyperlink display text on all
# three sheets, leaving hyperlink targets untouched.

$wb = $excel.ActiveWorkbook

$oldId = "4c3ad086-a9de-4630-8da0-a2acb5a7658a"
$newId = "bf2ba36c-0e8e-4eb8-8e64-e8dc319da171"
$oldHash = "328c4a9fa66d971c6e8210ba0abc9d999c7c06f7"
$newHash = "f7f65d20bcc6f852bb855cade74d66448e1a21ee"

$oldMdName = "$oldId.md"
$newMdName = "$newId.md"
$oldZhName = "$oldId.$oldHash.zh-cn.xlf"
$newZhName = "$newId.$newHash.zh-cn.xlf"
$oldDeName = "$oldId.$oldHash.de-de.xlf"
$newDeName = "$newId.$newHash.de-de.xlf"

$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/09536f2fac3ea5d2d9964074718f39a0a5d66424/e2e/$oldMdName"
$zhTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/538d110bad79e60d3637be12691ba93fad8fb370/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZhName"
$deTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fa9c387815ad9921ed2ad900fc8edf6aedd91fed/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDeName"

# ---- Overview sheet: just the handoff .md file name (A2) ----
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = $newMdName

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $mdTarget, "", "", $newMdName)

# ---- zh-cn sheet ----
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Value = $newMdName
$ws2.Range("D2").Value = $newZhName
$ws2.Range("E2").Value = "2016-03-13 19:04:58"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $mdTarget, "", "", $newMdName)
$ws2.Hyperlinks.Add($ws2.Range("B2"), $mdTarget, "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), $zhTarget, "", "", $newZhName)

# ---- de-de sheet ----
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Value = $newMdName
$ws3.Range("D2").Value = $newDeName
$ws3.Range("E2").Value = "2016-03-13 19:05:01"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $mdTarget, "", "", $newMdName)
$ws3.Hyperlinks.Add($ws3.Range("B2"), $mdTarget, "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), $deTarget, "", "", $newDeName)

"Done updating handoff report."
